$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 220.41667
$ws.Range("I11").Value = 220.41667
$ws.Range("K11").Value = 220.41667
$ws.Range("M11").Value = -80.41667000000001
$ws.Range("H15").Value = 451.96774
$ws.Range("I15").Value = 451.96774
$ws.Range("K15").Value = 1355.90322
$ws.Range("M15").Value = -1186.90322
$ws.Range("H40").Value = 13000
$ws.Range("J40").Value = 13000
$ws.Range("L40").Value = 13000
$ws.Range("N40").Value = -13350
$ws.Range("H86").Value = 1854.2858
$ws.Range("I86").Value = 996.3333
$ws.Range("J86").Value = 2497.75
$ws.Range("K86").Value = 996.3333
$ws.Range("L86").Value = 2497.75
$ws.Range("M86").Value = 126.6667
$ws.Range("N86").Value = -4743.75
$ws.Range("H89").Value = 1854.2858
$ws.Range("I89").Value = 996.3333
$ws.Range("J89").Value = 2497.75
$ws.Range("K89").Value = 4981.6665
$ws.Range("L89").Value = 12488.75
$ws.Range("M89").Value = 634.3334999999997
$ws.Range("N89").Value = -23720.75
$ws.Range("H116").Value = 812551.3
$ws.Range("I116").Value = 1348849.5
$ws.Range("J116").Value = 8104.1
$ws.Range("K116").Value = 1348849.5
$ws.Range("L116").Value = 8104.1
$ws.Range("M116").Value = -1345407.5
$ws.Range("N116").Value = -14988.1
$ws.Range("H138").Value = 1346577.1
$ws.Range("J138").Value = 2295435.8
$ws.Range("L138").Value = 6886307.399999999
$ws.Range("N138").Value = -6896587.399999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 20184.396
$ws.Range("I32").Value = 24542.318
$ws.Range("K32").Value = 24542.318
$ws.Range("M32").Value = -24255.318
$ws.Range("H102").Value = 1709.5834
$ws.Range("I102").Value = 1711.5
$ws.Range("K102").Value = 1711.5
$ws.Range("M102").Value = -89.5
$ws.Range("H110").Value = 22313.219
$ws.Range("I110").Value = 26649.879
$ws.Range("K110").Value = 26649.879
$ws.Range("M110").Value = -24604.879
$ws.Range("H132").Value = 1256.5217
$ws.Range("I132").Value = 1086.3636
$ws.Range("K132").Value = 3259.0908
$ws.Range("M132").Value = -729.0907999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1757.4193
$ws.Range("I86").Value = 1643.2222
$ws.Range("K86").Value = 1643.2222
$ws.Range("M86").Value = -520.2221999999999
$ws.Range("H89").Value = 1757.4193
$ws.Range("I89").Value = 1643.2222
$ws.Range("K89").Value = 8216.110999999999
$ws.Range("M89").Value = -2600.110999999999
$ws.Range("H105").Value = 2511.5557
$ws.Range("I105").Value = 2073.7932
$ws.Range("J105").Value = 4325.143
$ws.Range("K105").Value = 2073.7932
$ws.Range("L105").Value = 4325.143
$ws.Range("M105").Value = -326.7932000000001
$ws.Range("N105").Value = -7819.143
$ws.Range("H107").Value = 894.2083
$ws.Range("I107").Value = 850.6667
$ws.Range("K107").Value = 850.6667
$ws.Range("M107").Value = 1069.3333

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1726506
$ws.Range("J31").Value = 2740.4138
$ws.Range("L31").Value = 2740.4138
$ws.Range("N31").Value = -3330.4138
$ws.Range("H34").Value = 1726506
$ws.Range("J34").Value = 2740.4138
$ws.Range("L34").Value = 2740.4138
$ws.Range("N34").Value = -3144.4138
$ws.Range("H107").Value = 551.5
$ws.Range("I107").Value = 433.1
$ws.Range("J107").Value = 946.1667
$ws.Range("K107").Value = 433.1
$ws.Range("L107").Value = 946.1667
$ws.Range("M107").Value = 1486.9
$ws.Range("N107").Value = -4786.1667
$ws.Range("H132").Value = 36724.965
$ws.Range("I132").Value = 44997.22
$ws.Range("K132").Value = 134991.66
$ws.Range("M132").Value = -132461.66

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 346.77777
$ws.Range("I2").Value = 452.46155
$ws.Range("J2").Value = 72
$ws.Range("K2").Value = 2714.7693
$ws.Range("L2").Value = 432
$ws.Range("M2").Value = -2601.7693
$ws.Range("N2").Value = -658
$ws.Range("H37").Value = 42220.375
$ws.Range("J37").Value = 42220.375
$ws.Range("L37").Value = 126661.125
$ws.Range("N37").Value = -126885.125
$ws.Range("H68").Value = 2937.776
$ws.Range("J68").Value = 3173.04
$ws.Range("L68").Value = 9519.119999999999
$ws.Range("N68").Value = -11141.12
$ws.Range("H71").Value = 2937.776
$ws.Range("J71").Value = 3173.04
$ws.Range("L71").Value = 28557.36
$ws.Range("N71").Value = -36669.36
$ws.Range("H107").Value = 4963.9
$ws.Range("I107").Value = 1949.8572
$ws.Range("J107").Value = 5881.2173
$ws.Range("K107").Value = 5849.571599999999
$ws.Range("L107").Value = 17643.6519
$ws.Range("M107").Value = -3929.571599999999
$ws.Range("N107").Value = -21483.6519
$ws.Range("H131").Value = 2224.1052
$ws.Range("I131").Value = 3080.5625
$ws.Range("K131").Value = 9241.6875
$ws.Range("M131").Value = -4201.6875

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1734.35
$ws.Range("I113").Value = 826.1667
$ws.Range("K113").Value = 826.1667
$ws.Range("M113").Value = 1343.8333
$ws.Range("H126").Value = 5992
$ws.Range("J126").Value = 8011
$ws.Range("L126").Value = 24033
$ws.Range("N126").Value = -28973
$ws.Range("H132").Value = 2273.9395
$ws.Range("I132").Value = 1540.7307
$ws.Range("J132").Value = 4997.2856
$ws.Range("K132").Value = 4622.1921
$ws.Range("L132").Value = 14991.8568
$ws.Range("M132").Value = -2092.1921
$ws.Range("N132").Value = -20051.8568

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2625.52
$ws.Range("I40").Value = 2460.9048
$ws.Range("K40").Value = 2460.9048
$ws.Range("M40").Value = -2324.9048
$ws.Range("H46").Value = 4608.2
$ws.Range("J46").Value = 4534.9473
$ws.Range("L46").Value = 4534.9473
$ws.Range("N46").Value = -4910.9473
$ws.Range("H82").Value = 3123.7097
$ws.Range("I82").Value = 2956
$ws.Range("J82").Value = 3327.3572
$ws.Range("K82").Value = 2956
$ws.Range("L82").Value = 3327.3572
$ws.Range("M82").Value = -2595
$ws.Range("N82").Value = -4049.3572
$ws.Range("H85").Value = 3123.7097
$ws.Range("I85").Value = 2956
$ws.Range("J85").Value = 3327.3572
$ws.Range("K85").Value = 2956
$ws.Range("L85").Value = 3327.3572
$ws.Range("M85").Value = -1708
$ws.Range("N85").Value = -5823.3572
$ws.Range("H136").Value = 4534.364
$ws.Range("I136").Value = 3250.8125
$ws.Range("J136").Value = 7957.1665
$ws.Range("K136").Value = 9752.4375
$ws.Range("L136").Value = 23871.4995
$ws.Range("M136").Value = -7202.4375
$ws.Range("N136").Value = -28971.4995

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H23").Value = 254
$ws.Range("I23").Value = 102.5
$ws.Range("J23").Value = 405.5
$ws.Range("K23").Value = 102.5
$ws.Range("L23").Value = 405.5
$ws.Range("M23").Value = 126.5
$ws.Range("N23").Value = -863.5
$ws.Range("H54").Value = 14833.333
$ws.Range("H122").Value = 31389.043
$ws.Range("I122").Value = 47522.965
$ws.Range("K122").Value = 142568.895
$ws.Range("M122").Value = -140118.895
$ws.Range("H132").Value = 24822.576
$ws.Range("I132").Value = 31341.8
$ws.Range("J132").Value = 4450
$ws.Range("K132").Value = 94025.39999999999
$ws.Range("L132").Value = 13350
$ws.Range("M132").Value = -91495.39999999999
$ws.Range("N132").Value = -18410
$ws.Range("H135").Value = 123290.336
$ws.Range("J135").Value = 123290.336
$ws.Range("L135").Value = 123290.336
$ws.Range("N135").Value = -133430.336
$ws.Range("H136").Value = 18494.066
$ws.Range("I136").Value = 19480.875
$ws.Range("K136").Value = 58442.625
$ws.Range("M136").Value = -55892.625
